# Timeplan.xlsx edit:
#  - Insert a new milestone row ("Override attack methods for races",
#    Estimated 60, Real 85) right above the "Sum" row.
#  - The "Sum" row (now shifted down one row) keeps summing the full
#    milestone range, now through the newly-inserted row.
#  - Leave the final hidden helper row intact (it just shifts down too).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push the "Sum" row (and the hidden row below it) down by one to make
# room for the new milestone entry.
$ws.Rows("7:7").Insert()

# Match the formatting of the existing milestone rows (e.g. row 3) for
# the newly inserted row.
$ws.Range("A3").Copy()
$ws.Range("A7").PasteSpecial(-4122)
$ws.Range("B3:C3").Copy()
$ws.Range("B7:C7").PasteSpecial(-4122)

# New milestone data.
$ws.Range("A7").Value = "Override attack methods for races"
$ws.Range("B7").Value = 60
$ws.Range("C7").Value = 85

# Refresh the Sum row's formulas to include the new row.
$ws.Range("B8").Formula = "=SUM(B2:B7)"
$ws.Range("C8").Formula = "=SUM(C2:C7)"

# Match the author's final selection.
$ws.Range("A8").Select()
